$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new row 19 data ---
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = 2.0937847222222223
$ws.Range("B19").NumberFormat = $ws.Range("B18").NumberFormat
$ws.Range("C19").Formula = "=SUM(B2:B19)+1.2708333333"
$ws.Range("C19").NumberFormat = $ws.Range("C18").NumberFormat
$ws.Range("D19").Value = "El desorden que dejas (Audiovisual, Spanish, New):37; Perdida 2018 (Subtitled, Spanish, New):40; Siempre bruja (Audiovisual, Spanish, New):39; Corazón loco (Audiovisual, Spanish, New):33; [Este Jugador Tiene Un Survival HERMOSO en la Nether Update]https://youtu.be/lVYn6fshxFY) (Audiovisual, Spanish, New):32; Amadeus (Audiovisual, English, Familiar):37;"

# --- Move the selection the way it ended up after the edit ---
$ws.Range("C20").Select()
